$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.153.42"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.779.20"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.26"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.68"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "2.036.20"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.93"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.774.93"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "34.109.37"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.85"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.35"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.62"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  +4.09%  "
$ws.Range("E33").Value = "  +5.84%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "1.444.20"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.660"
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("E37").Value = "  +6.35%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.19"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.44"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "1.938.20"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.01"
$ws.Range("E51").Value = "  -1.45%  "
